# Indiana roster update: swap in "James Johnson" (replacing the TW/Exhibit-10
# ordering), shift the trailing veteran rows down, and drop "Serge Ibaka"
# entirely. Net effect: row 13 becomes James Johnson, rows 14/16/17/18 shift
# to host the players that used to occupy the row below/above them, and
# Daniel Theis (row 15) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Trevelin Queen (TW) -> James Johnson
$ws.Range("B13").Value = 16
$ws.Range("C13").Value = "James Johnson"
$ws.Range("D13").Value = "PF"
$ws.Range("E13").Value = "6-7"
$ws.Range("F13").Value = 240
$ws.Range("G13").Value = "February 20, 1987"
$ws.Range("I13").Value = "'13"
$ws.Range("J13").Value = "Wake Forest"
$ws.Range("K13").Value = "https://www.basketball-reference.com/players/j/johnsja01.html"

# Row 14: Kendall Brown (TW) -> Trevelin Queen (TW)
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = "Trevelin Queen (TW)"
$ws.Range("E14").Value = "6-6"
$ws.Range("F14").Value = 190
$ws.Range("G14").Value = "February 25, 1997"
$ws.Range("I14").Value = "'1"
$ws.Range("J14").Value = "College of Marin, New Mexico Military Institute, New Mexico State"
$ws.Range("K14").Value = "https://www.basketball-reference.com/players/q/queentr01.html"

# Row 15: Daniel Theis stays Daniel Theis (values unchanged)

# Row 16: George Hill -> Kendall Brown (TW)
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "Kendall Brown (TW)"
$ws.Range("D16").Value = "SG"
$ws.Range("E16").Value = "6-8"
$ws.Range("F16").Value = 205
$ws.Range("G16").Value = "May 11, 2003"
$ws.Range("I16").Value = "R"
$ws.Range("J16").Value = "Baylor"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/b/brownke03.html"

# Row 17: Serge Ibaka -> Jordan Nwora
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = "Jordan Nwora"
$ws.Range("D17").Value = "SF"
$ws.Range("E17").Value = "6-8"
$ws.Range("F17").Value = 225
$ws.Range("G17").Value = "September 9, 1998"
$ws.Range("H17").Value = "us"
$ws.Range("I17").Value = "'2"
$ws.Range("J17").Value = "Louisville"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/n/nworajo01.html"

# Row 18: Jordan Nwora -> George Hill
$ws.Range("B18").Value = 7
$ws.Range("C18").Value = "George Hill"
$ws.Range("D18").Value = "PG"
$ws.Range("E18").Value = "6-4"
$ws.Range("F18").Value = 188
$ws.Range("G18").Value = "May 4, 1986"
$ws.Range("I18").Value = "'14"
$ws.Range("J18").Value = "IUPUI"
$ws.Range("K18").Value = "https://www.basketball-reference.com/players/h/hillge01.html"
